# Update the ObjTables/SBtab header metadata (date + objTablesVersion) that is
# stamped into cell A1 (and, for the very first sheet, also A2) of every
# worksheet in the workbook.
#
#   objTablesVersion: '0.0.9' -> '1.0.0'
#   date             : '2020-04-27 01:09:00'/'...01' -> '2020-05-29 00:22:55'

$wb = $excel.ActiveWorkbook

$oldDate = "2020-04-27 01:09:0"   # shared prefix; last digit varies (0 or 1) across sheets
$newDate = "2020-05-29 00:22:55"
$oldVersion = "objTablesVersion='0.0.9'"
$newVersion = "objTablesVersion='1.0.0'"

function Update-ObjTablesCell($range) {
    $text = $range.Value()
    if ($text -eq $null) { return }
    if ($text -like "*$oldDate*") {
        # Replace date='2020-04-27 01:09:0X' -> date='2020-05-29 00:22:55'
        $text = [System.Text.RegularExpressions.Regex]::Replace($text, "date='2020-04-27 01:09:0[01]'", "date='$newDate'")
        $text = $text.Replace($oldVersion, $newVersion)
        $range.Value = $text
    }
}

# The classes in workbook/tab order -- each one (except the first) owns a
# worksheet named "!!<Class>" whose A1 cell holds the per-table ObjTables
# header comment.
$classes = @(
    "Compartment",
    "Compound",
    "Definition",
    "Enzyme",
    "FbcObjective",
    "Gene",
    "Layout",
    "Measurement",
    "PbConfig",
    "Position",
    "Protein",
    "Quantity",
    "QuantityInfo",
    "QuantityMatrix",
    "Reaction",
    "ReactionStoichiometry",
    "Regulator",
    "Relation",
    "Relationship",
    "SparseMatrix",
    "SparseMatrixColumn",
    "SparseMatrixOrdered",
    "SparseMatrixRow",
    "StoichiometricMatrix",
    "rxnconContingencyList",
    "rxnconReactionList"
)

$first = $true
foreach ($class in $classes) {
    $ws = $wb.Worksheets.Item("!!$class")

    # Sheets ship protected (no password) so cell writes must unprotect first.
    $ws.Unprotect()

    # Top-left cell: per-table ObjTables header ('!!ObjTables ... class=... date=... objTablesVersion=...')
    Update-ObjTablesCell $ws.Range("A1")

    if ($first) {
        # The very first worksheet's A2 carries the file-level header
        # ('!!!ObjTables schema=... objTablesVersion=... date=...') as well
        # as its own class header in A1 handled above.
        Update-ObjTablesCell $ws.Range("A2")
        $first = $false
    }

    # Restore sheet protection to match the original document state.
    $ws.Protect()
}
